$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.836
$ws.Range("C2").Value = 0.915
$ws.Range("D2").Value = 0.964
$ws.Range("E2").Value = 0.836
$ws.Range("F2").Value = 0.53
$ws.Range("G2").Value = 0.832

# Row 3
$ws.Range("B3").Value = 0.903
$ws.Range("C3").Value = 0.928
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.836
$ws.Range("F3").Value = 0.522
$ws.Range("G3").Value = 0.842
